$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Docentes responsáveis:" value row (old row 13, containing only
# B13/C13 = "5840650 - Janaína Ferreira Batista") is removed entirely;
# everything below shifts up by one row.
$ws.Rows("13:13").Delete()

# After the shift, a handful of cells get new content (values that were
# effectively "promoted" from other rows / replaced outright).
$ws.Range("B10").Value = "5840650 - Janaína Ferreira Batista"
$ws.Range("C10").Value = "5840650 - Janaína Ferreira Batista"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "5840650 - Janaína Ferreira Batista"
$ws.Range("C18").Value = "5840650 - Janaína Ferreira Batista"

$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
